$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 18; $row++) {
    $cell = $ws.Cells.Item($row, 5)  # Column E = 5
    if ($cell.Value() -eq "fullRNASEQ") {
        $cell.Value = "fullRNASeq"
    }
}
